{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block (and\n// the blank paragraph that separates it from the bibliography) that was\n// scraped from the Jekyll site footer.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst JUPITER_TEXT = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst COPYRIGHT_TEXT =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n// Locate the \"Ver no Jupiter ...\" paragraph.\nlet jupiterIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === JUPITER_TEXT) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nif (jupiterIndex !== -1) {\n  // The paragraph right after it should be the copyright/footer line.\n  const copyrightIndex = jupiterIndex + 1;\n  // The paragraph right before it is the blank separator paragraph that\n  // was inserted right after the bibliography's last entry.\n  const blankIndex = jupiterIndex - 1;\n\n  if (copyrightIndex < items.length && items[copyrightIndex].text === COPYRIGHT_TEXT) {\n    items[copyrightIndex].delete();\n  }\n\n  items[jupiterIndex].delete();\n\n  if (blankIndex >= 0 && items[blankIndex].text === \"\") {\n    items[blankIndex].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph separating it from the bibliography) that was\n# scraped from the Jekyll site footer.\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$JUPITER_TEXT = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$COPYRIGHT_TEXT = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n# Locate the \"Ver no Jupiter ...\" paragraph (strip the trailing paragraph\n# mark / cell mark before comparing).\n$jupiterIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $paras.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $JUPITER_TEXT) {\n        $jupiterIndex = $i\n        break\n    }\n}\n\nif ($jupiterIndex -ne -1) {\n    $copyrightIndex = $jupiterIndex + 1\n    $blankIndex = $jupiterIndex - 1\n\n    # Delete from the bottom up so earlier indices stay valid.\n    if ($copyrightIndex -le $paras.Count) {\n        $copyrightText = $paras.Item($copyrightIndex).Range.Text.TrimEnd([char]13, [char]7)\n        if ($copyrightText -eq $COPYRIGHT_TEXT) {\n            $paras.Item($copyrightIndex).Range.Delete()\n        }\n    }\n\n    $paras.Item($jupiterIndex).Range.Delete()\n\n    if ($blankIndex -ge 1) {\n        $blankText = $paras.Item($blankIndex).Range.Text.TrimEnd([char]13, [char]7)\n        if ($blankText -eq \"\") {\n            $paras.Item($blankIndex).Range.Delete()\n        }\n    }\n}\n"}
